# Commit: "Added OSATS results and removed grasp results for now."
#
# This script:
#  1. Populates the previously-empty "scale_OSATS" worksheet with a header
#     row plus 7 data rows (new OSATS meta-analysis results), including the
#     derived-statistic formulas (O:R) and a couple of per-row helper
#     formulas in columns K/N that reconstruct SD values from ranges.
#  2. Updates view/selection state: "tool_movements" loses the selected-tab
#     flag, "pupil_blinks" gets a new selection (H29), and "scale_OSATS"
#     becomes the selected tab with selection M10 (matching the new
#     workbookView.activeTab pointing at the last sheet).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) Update "tool_movements": it is no longer the tab shown when the
#    workbook opens (selection itself is unchanged).
# ---------------------------------------------------------------------
$wsMovements = $wb.Worksheets.Item("tool_movements")
$wsMovements.Range("G12").Select()

# ---------------------------------------------------------------------
# 2) Update "pupil_blinks": selection moves from A1:R1 to H29.
# ---------------------------------------------------------------------
$wsBlinks = $wb.Worksheets.Item("pupil_blinks")
$wsBlinks.Range("H29").Select()

# ---------------------------------------------------------------------
# 3) Populate "scale_OSATS" with the new OSATS data table.
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("scale_OSATS")

# --- Header row -------------------------------------------------------
$ws.Range("A1").Value = "i"
$ws.Range("B1").Value = "Author"
$ws.Range("C1").Value = "Year"
$ws.Range("D1").Value = "Study"
$ws.Range("E1").Value = "Journal"
$ws.Range("F1").Value = "Technique"
$ws.Range("G1").Value = "Task"
$ws.Range("H1").Value = "Note"
$ws.Range("I1").Value = "Nn"
$ws.Range("J1").Value = "Mn"
$ws.Range("K1").Value = "SDn"
$ws.Range("L1").Value = "Ne"
$ws.Range("M1").Value = "Me"
$ws.Range("N1").Value = "SDe"
$ws.Range("O1").Value = "SDpooled"
$ws.Range("P1").Value = "SMD"
$ws.Range("Q1").Value = "g"
$ws.Range("R1").Value = "SDg"

# --- Row 2: Nickel et al. ---------------------------------------------
$ws.Range("A2").Value = 0
$ws.Range("B2").Value = "Nickel et al."
$ws.Range("C2").Value = 2016
$ws.Range("D2").Value = "Direct Observation versus Endoscopic Video Recording-Based Rating with the Objective Structured Assessment of Technical Skills for Training of Laparoscopic Cholecystectomy"
$ws.Range("E2").Value = "European Surgical Research"
$ws.Range("F2").Value = "Laparoscopy"
$ws.Range("H2").Value = "OSATS score from Table 1, direct observation, novices and experts compared"
$ws.Range("I2").Value = 15
$ws.Range("J2").Value = 48.6
$ws.Range("K2").Value = 8.8000000000000007
$ws.Range("L2").Value = 6
$ws.Range("M2").Value = 65.3
$ws.Range("N2").Value = 10.4

# --- Row 3: Paley et al. -----------------------------------------------
$ws.Range("A3").Value = 1
$ws.Range("B3").Value = "Paley et al."
$ws.Range("C3").Value = 2021
$ws.Range("D3").Value = "Crowdsourced Assessment of Surgical Skill Proficiency in Cataract Surgery"
$ws.Range("E3").Value = "Journal of Surgical Education"
$ws.Range("H3").Value = "Used modified OSATS. SD estimated from Figure 1F. Used expert ratings."
$ws.Range("I3").Value = 6
$ws.Range("J3").Value = 7.3
$ws.Range("K3").Value = 2
$ws.Range("L3").Value = 9
$ws.Range("M3").Value = 21
$ws.Range("N3").Value = 4

# --- Row 4: Kassab et al. ----------------------------------------------
$ws.Range("A4").Value = 2
$ws.Range("B4").Value = "Kassab et al."
$ws.Range("C4").Value = 2011
$ws.Range("D4").Value = '"Blowing up the barriers" in surgical training: Exploring and validating the concept of distributed simulation'
$ws.Range("E4").Value = "Annals of Surgery"
$ws.Range("F4").Value = "Laparoscopy"
$ws.Range("G4").Value = "Box trainer"
$ws.Range("H4").Value = "Study had two tasks, results are for DS (distributed simulation) because these results were given in the text (box trainer results only as figure). Note that DS was novel task developed for this study."
$ws.Range("I4").Value = 10
$ws.Range("J4").Value = 16.3
$ws.Range("K4").Value = 3.8
$ws.Range("L4").Value = 10
$ws.Range("M4").Value = 27.3
$ws.Range("N4").Value = 5.7

# --- Row 5: Black et al. -----------------------------------------------
$ws.Range("A5").Value = 3
$ws.Range("B5").Value = "Black et al."
$ws.Range("C5").Value = 2010
$ws.Range("D5").Value = "Assessment of surgical competence at carotid endarterectomy under local anaesthesia in a simulated operating theatre"
$ws.Range("E5").Value = "British Journal of Surgery"
$ws.Range("H5").Value = "Results for crisis scenario"
$ws.Range("I5").Value = 10
$ws.Range("J5").Value = 15.5
$ws.Range("K5").Formula = "=(19-12)*(3/4)"
$ws.Range("L5").Value = 10
$ws.Range("M5").Value = 36
$ws.Range("N5").Formula = "=(36-35)*(3/4)"

# --- Row 6: Willems et al. ---------------------------------------------
$ws.Range("A6").Value = 4
$ws.Range("B6").Value = "Willems et al."
$ws.Range("C6").Value = 2009
$ws.Range("D6").Value = "Assessing Endovascular Skills using the Simulator for Testing and Rating Endovascular Skills (STRESS) Machine"
$ws.Range("E6").Value = "European Journal of Vascular and Endovascular Surgery"
$ws.Range("H6").Value = "Combination of OSATS and some other score? May not be suitable for comparison here. Remove in the future. SDs estimated from Figure 2."
$ws.Range("I6").Value = 8
$ws.Range("J6").Value = 42.75
$ws.Range("K6").Value = 5
$ws.Range("L6").Value = 5
$ws.Range("M6").Value = 82.8
$ws.Range("N6").Formula = "=25*(3/5)"

# --- Row 7: Leong et al. -----------------------------------------------
$ws.Range("A7").Value = 5
$ws.Range("B7").Value = "Leong et al."
$ws.Range("C7").Value = 2008
$ws.Range("D7").Value = "Validation of orthopaedic bench models for trauma surgery"
$ws.Range("E7").Value = "Journal of Bone and Joint Surgery - Series B"
$ws.Range("H7").Value = "Used results for DCP, dynamic comperssion plate. Esimtaed values from boxplot."
$ws.Range("I7").Value = 8
$ws.Range("J7").Value = 35
$ws.Range("K7").Formula = "=15*(3/4)"
$ws.Range("L7").Value = 6
$ws.Range("M7").Value = 65
$ws.Range("N7").Formula = "=12*(3/4)"

# --- Row 8: Hance et al. -----------------------------------------------
$ws.Range("A8").Value = 6
$ws.Range("B8").Value = "Hance et al."
$ws.Range("C8").Value = 2005
$ws.Range("D8").Value = "Objective assessment of technical skills in cardiac surgery"
$ws.Range("E8").Value = "European Journal of Cardio-thoracic Surgery"
$ws.Range("G8").Value = "LAD anastomosis"
$ws.Range("H8").Value = "Paper reported several tasks, live and blinded scoring. Values here are for LAD anastomosis, blinded scoring."
$ws.Range("I8").Value = 12
$ws.Range("J8").Value = 15.5
$ws.Range("K8").Formula = "=(19.5-13.25)*(3/4)"
$ws.Range("L8").Value = 13
$ws.Range("M8").Value = 24
$ws.Range("N8").Formula = "=(34-21)*(3/4)"

# --- Derived-statistic formulas for every data row (O:R) ---------------
for ($r = 2; $r -le 8; $r++) {
    $ws.Range("O$r").Formula = "=SQRT(((I$r-1)*POWER(K$r,2) + (L$r-1)*POWER(N$r,2))/((I$r-1)+(L$r-1)))"
    $ws.Range("P$r").Formula = "=(J$r-M$r)/O$r"
    $ws.Range("Q$r").Formula = "=P$r*(1- (3/(4*(I$r+L$r)-9)))"
    $ws.Range("R$r").Formula = "=SQRT((I$r+L$r)/(I$r*L$r)+(POWER(P$r,2)/(2*(I$r+L$r))))"
}

# --- Final selection / active-tab state --------------------------------
$ws.Range("M10").Select()
